# Changed x axis on sample brightness graph
# Updates crop-begin / crop-end time values for a handful of rows, re-applies
# the "Arial" cell format to the data rows that were still carrying the old
# default format, and flags row 29 as Questionable (like rows 12/20/28/31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: pure value edits (keep existing cell formatting) ---------
$ws.Range("F2").Value2 = 9.0
$ws.Range("G2").Value2 = 18.5

$ws.Range("F3").Value2 = 17.0
$ws.Range("G3").Value2 = 42.0

$ws.Range("F32").Value2 = 6.5
$ws.Range("G32").Value2 = 18.5

$ws.Range("F33").Value2 = 10.0
$ws.Range("G33").Value2 = 41.0

# --- Phase 2: re-format the remaining "default" styled cells to match ---
# the rest of the table (copy the already-correct format from E4 which
# uses the normal Arial / theme text style).
$ws.Range("E4").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33:E36").PasteSpecial(-4122)
$ws.Range("E7:E31").PasteSpecial(-4122)
$ws.Range("F4:G31").PasteSpecial(-4122)
$ws.Range("F34:G36").PasteSpecial(-4122)
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H31").PasteSpecial(-4122)

# --- Phase 3: new "Questionable" flag on row 29 --------------------------
$ws.Range("H29").Value2 = "Yes"
$ws.Range("H29").Font.ThemeColor = 1
$ws.Range("H29").Font.Name = "Serif"

$excel.CutCopyMode = 0
